# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) YDS sheet - append the new game's per-play yardage log entries to the
#    four running logs (Rush/Pass, Offense/Defense).
# ---------------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Text + " 6 3 9 0 3 -3 3 0 4 5 4 25 -2 -1 1 6 5 29 -4 4 8"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Text + " 3 7 3 25 40 -3 44 9 4 19 7 45 26 21"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Text + " 5 4 6 0 0 3 11 1 3 1 8 6 4 9 3 4 3 3 -1 7 1 0 4 1"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Text + " 14 16 3 5 5 19 11 16 6 18 -1 14 5 8 23 8 20 13 19 6 4 4 6 18 7 7 6 3 11 2 10 6"

# ---------------------------------------------------------------------------
# 2) OFF sheet - updated season totals after the new game.
# ---------------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("C2").Value = 422
$offWs.Range("F2").Value = 135
$offWs.Range("G2").Value = 137
$offWs.Range("H2").Value = 11
$offWs.Range("J2").Value = 57
$offWs.Range("L2").Value = 537
$offWs.Range("M2").Value = 361
$offWs.Range("Q2").Value = 1007
$offWs.Range("B3").Value = 22
$offWs.Range("C3").Value = 351
$offWs.Range("E3").Value = 68
$offWs.Range("H3").Value = 67
$offWs.Range("I3").Value = 126
$offWs.Range("J3").Value = 99
$offWs.Range("N3").Value = 31

# ---------------------------------------------------------------------------
# 3) DEF sheet - updated season totals after the new game.
# ---------------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 455
$defWs.Range("E2").Value = 19
$defWs.Range("F2").Value = 136
$defWs.Range("G2").Value = 139
$defWs.Range("I2").Value = 13
$defWs.Range("J2").Value = 65
$defWs.Range("L2").Value = 586
$defWs.Range("M2").Value = 387
$defWs.Range("O2").Value = 54
$defWs.Range("P2").Value = 31
$defWs.Range("Q2").Value = 1083
$defWs.Range("C3").Value = 368
$defWs.Range("E3").Value = 44
$defWs.Range("F3").Value = 213
$defWs.Range("G3").Value = 96
$defWs.Range("H3").Value = 45
$defWs.Range("I3").Value = 115
$defWs.Range("J3").Value = 106
$defWs.Range("N3").Value = 36

# ---------------------------------------------------------------------------
# 4) ST sheet - updated special-teams totals + appended per-kick logs.
# ---------------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = 173
$stWs.Range("D2").Value = 130
$stWs.Range("F2").Value = 78
$stWs.Range("G2").Value = 70
$stWs.Range("J2").Value = 69
$stWs.Range("K2").Value = 58
$stWs.Range("B3").Value = 130

$stWs.Range("B4").Value = $stWs.Range("B4").Text + " 64"
$stWs.Range("B5").Value = $stWs.Range("B5").Text + " 36"
$stWs.Range("B6").Value = $stWs.Range("B6").Text + " 14 24"
$stWs.Range("D3").Value = $stWs.Range("D3").Text + " 42 29 52 57 52 40"
$stWs.Range("D4").Value = $stWs.Range("D4").Text + " 17 0 9 11 10 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Text + " 0"

# ---------------------------------------------------------------------------
# 5) TURNS sheet - updated turnover totals.
# ---------------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("C2").Value = 16
$turnsWs.Range("D2").Value = 16

# ---------------------------------------------------------------------------
# 6) PEN sheet - updated penalty totals.
# ---------------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("D2").Value = 12
$penWs.Range("B3").Value = 50
$penWs.Range("D4").Value = 16

Write-Output "Applied Wild Card round simulation results"
